$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 24
$ws.Cells.Item(2, 2).Value = " fwd header length"
$ws.Cells.Item(2, 3).Value = 0.05352177463688457
$ws.Cells.Item(3, 1).Value = 25
$ws.Cells.Item(3, 2).Value = " bwd header length"
$ws.Cells.Item(3, 3).Value = 0.05132114164522052
$ws.Cells.Item(4, 1).Value = 39
$ws.Cells.Item(4, 2).Value = " idle max"
$ws.Cells.Item(4, 3).Value = 0.05079032304311388
$ws.Cells.Item(5, 1).Value = 26
$ws.Cells.Item(5, 2).Value = "fwd packets/s"
$ws.Cells.Item(5, 3).Value = 0.04765685279244385
$ws.Cells.Item(6, 1).Value = 20
$ws.Cells.Item(6, 2).Value = " bwd iat mean"
$ws.Cells.Item(6, 3).Value = 0.04715369417577813
$ws.Cells.Item(7, 1).Value = 19
$ws.Cells.Item(7, 2).Value = "bwd iat total"
$ws.Cells.Item(7, 3).Value = 0.04525872297014155
$ws.Cells.Item(8, 1).Value = 37
$ws.Cells.Item(8, 2).Value = "idle mean"
$ws.Cells.Item(8, 3).Value = 0.0398545026066294
$ws.Cells.Item(9, 1).Value = 21
$ws.Cells.Item(9, 2).Value = " bwd iat std"
$ws.Cells.Item(9, 3).Value = 0.03650224955751528
$ws.Cells.Item(10, 1).Value = 18
$ws.Cells.Item(10, 2).Value = " fwd iat min"
$ws.Cells.Item(10, 3).Value = 0.03148568619238207
$ws.Cells.Item(11, 1).Value = 23
$ws.Cells.Item(11, 2).Value = " bwd iat min"
$ws.Cells.Item(11, 3).Value = 0.02858958990602306
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = " flow duration"
$ws.Cells.Item(12, 3).Value = 0.02721170215440187
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = "unnamed: 0"
$ws.Cells.Item(13, 3).Value = 0.02497444816996061
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = " source port"
$ws.Cells.Item(14, 3).Value = 0.02344772370138886
$ws.Cells.Item(15, 1).Value = 2
$ws.Cells.Item(15, 2).Value = " destination port"
$ws.Cells.Item(15, 3).Value = 0.02326241074107733
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = " total fwd packets"
$ws.Cells.Item(16, 3).Value = 0.01619763551876954
$ws.Cells.Item(17, 1).Value = 22
$ws.Cells.Item(17, 2).Value = " bwd iat max"
$ws.Cells.Item(17, 3).Value = 0.00868639303857755
$ws.Cells.Item(18, 1).Value = 32
$ws.Cells.Item(18, 2).Value = " init_win_bytes_backward"
$ws.Cells.Item(18, 3).Value = 0.007685638186746611
$ws.Cells.Item(19, 1).Value = 28
$ws.Cells.Item(19, 2).Value = " fwd header length.1"
$ws.Cells.Item(19, 3).Value = 0.007272873338239656
$ws.Cells.Item(20, 1).Value = 31
$ws.Cells.Item(20, 2).Value = "init_win_bytes_forward"
$ws.Cells.Item(20, 3).Value = 0.007086391199984689
$ws.Cells.Item(21, 1).Value = 29
$ws.Cells.Item(21, 2).Value = "subflow fwd packets"
$ws.Cells.Item(21, 3).Value = 0.006994157484101189
$ws.Cells.Item(22, 1).Value = 27
$ws.Cells.Item(22, 2).Value = " bwd packets/s"
$ws.Cells.Item(22, 3).Value = 0.006187803114247391
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = " total backward packets"
$ws.Cells.Item(23, 3).Value = 0.00578090648389456
$ws.Cells.Item(24, 1).Value = 30
$ws.Cells.Item(24, 2).Value = " subflow bwd packets"
$ws.Cells.Item(24, 3).Value = 0.005058150279394818
$ws.Cells.Item(25, 1).Value = 38
$ws.Cells.Item(25, 2).Value = " idle std"
$ws.Cells.Item(25, 3).Value = 0.004314152961778288
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = " protocol"
$ws.Cells.Item(26, 3).Value = 0.003260901756261592
$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = " fwd packet length max"
$ws.Cells.Item(27, 3).Value = 0.003239827766553275
$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = " flow packets/s"
$ws.Cells.Item(28, 3).Value = 0.001769433015466189
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = " flow iat mean"
$ws.Cells.Item(29, 3).Value = 0.001540208623576104
$ws.Cells.Item(30, 1).Value = 35
$ws.Cells.Item(30, 2).Value = " active max"
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = 36
$ws.Cells.Item(31, 2).Value = " active min"
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 1).Value = 34
$ws.Cells.Item(32, 2).Value = " active std"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 1).Value = 8
$ws.Cells.Item(33, 2).Value = " bwd packet length mean"
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 1).Value = 12
$ws.Cells.Item(34, 2).Value = " flow iat max"
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 1).Value = 13
$ws.Cells.Item(35, 2).Value = " flow iat min"
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = " flow iat std"
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 1).Value = 17
$ws.Cells.Item(37, 2).Value = " fwd iat max"
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(38, 1).Value = 15
$ws.Cells.Item(38, 2).Value = " fwd iat mean"
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 1).Value = 16
$ws.Cells.Item(39, 2).Value = " fwd iat std"
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = 33
$ws.Cells.Item(40, 2).Value = "active mean"
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 1).Value = 14
$ws.Cells.Item(41, 2).Value = "fwd iat total"
$ws.Cells.Item(41, 3).Value = 0

Write-Host "applied"